$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the values from columns O,R,AN,AQ (rows 1-3) into B,C,D,E (rows 1-3)
$ws.Range("O1:O3").Copy() | Out-Null
$ws.Range("B1:B3").PasteSpecial(-4163) | Out-Null

$ws.Range("R1:R3").Copy() | Out-Null
$ws.Range("C1:C3").PasteSpecial(-4163) | Out-Null

$ws.Range("AN1:AN3").Copy() | Out-Null
$ws.Range("D1:D3").PasteSpecial(-4163) | Out-Null

$ws.Range("AQ1:AQ3").Copy() | Out-Null
$ws.Range("E1:E3").PasteSpecial(-4163) | Out-Null

$ws.Range("B1:E3").Select() | Out-Null
